$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (IT / Daycare): cells C,D,E,G,H change from "X" to "XX" ---
$ws.Range("C6").Value = "XX"
$ws.Range("D6").Value = "XX"
$ws.Range("E6").Value = "XX"
$ws.Range("G6").Value = "XX"
$ws.Range("H6").Value = "XX"
# F6 keeps "X" but is highlighted in red (new font: Calibri, red)
$ws.Range("F6").Font.Color = 255

# --- Row 7 (IT / Home): cells C,D,E,G,H change from "X" to "XX" ---
$ws.Range("C7").Value = "XX"
$ws.Range("D7").Value = "XX"
$ws.Range("E7").Value = "XX"
$ws.Range("G7").Value = "XX"
$ws.Range("H7").Value = "XX"
# F7 keeps "X" but is highlighted in red (reuses existing red "Calibri (Body)" font)
$f7 = $ws.Range("F7")
$f7.Font.Color = 255
$f7.Font.Name = "Calibri (Body)"

# --- Legend text update (row 15) ---
$ws.Range("A15").Value = "XX = integrated into norms output"

# --- Update active cell selection ---
$ws.Range("H7").Select()
